# Refresh the cryptos price list ("Updated cryptos list ... with GitHub
# Actions"): update Price/Volume(1h) figures for most rows, and re-sort the
# Stacks / LidoDAOToken / FirstDigitalUSD trio (rows 39-41) into their new
# order with the values that go with each coin.
#
# Price-column cells ('D') hold values like "394.18" or "1.00" that Excel's
# own type-inference would otherwise coerce into real numbers. Forcing
# NumberFormat to Text ("@") before writing keeps them as plain text, which
# matches how the source file stores these cells (inline strings, not
# numeric values).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '57.107.79'
$ws.Range('E2').Value = '  +7.09%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.237.82'
$ws.Range('E3').Value = '  +2.77%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '394.18'
$ws.Range('E5').Value = '  -1.03%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '106.61'
$ws.Range('E6').Value = '  -0.26%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.234.66'
$ws.Range('E7').Value = '  +2.84%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.564'
$ws.Range('E8').Value = '  +3.48%  '
$ws.Range('E9').Value = '  -0.01%  '
$ws.Range('E10').Value = '  +1.31%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '38.84'
$ws.Range('E11').Value = '  -0.35%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0973'
$ws.Range('E12').Value = '  +11.83%  '
$ws.Range('E13').Value = '  +1.69%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.750.40'
$ws.Range('E14').Value = '  +2.91%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '8.12'
$ws.Range('E15').Value = '  +1.82%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '18.89'
$ws.Range('E16').Value = '  -0.36%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.242.10'
$ws.Range('E17').Value = '  +3.03%  '
$ws.Range('E18').Value = '  -1.91%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.01'
$ws.Range('E19').Value = '  +2.21%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '56.863.44'
$ws.Range('E20').Value = '  +6.80%  '
$ws.Range('E21').Value = '  +1.36%  '
$ws.Range('E22').Value = '  +7.79%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '12.91'
$ws.Range('E23').Value = '  +0.14%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '298.17'
$ws.Range('E24').Value = '  +10.03%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '73.73'
$ws.Range('E25').Value = '  +3.58%  '
$ws.Range('E26').Value = '  -3.01%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '4.39'
$ws.Range('E27').Value = '  +3.50%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '27.82'
$ws.Range('E28').Value = '  +0.69%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.72'
$ws.Range('E29').Value = '  -4.03%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.22'
$ws.Range('E30').Value = '  -3.82%  '
$ws.Range('E31').Value = '  -1.79%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.00'
$ws.Range('E32').Value = '  +0.00%  '
$ws.Range('E33').Value = '  -1.14%  '
$ws.Range('E34').Value = '  -0.86%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '37.29'
$ws.Range('E35').Value = '  -0.24%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0482'
$ws.Range('E36').Value = '  -1.99%  '
$ws.Range('E37').Value = '  +1.00%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '51.69'
$ws.Range('E38').Value = '  +2.41%  '
$ws.Range('B39').Value = 'LidoDAOToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.51'
$ws.Range('E39').Value = '  -1.56%  '
$ws.Range('B40').Value = 'FirstDigitalUSD'
$ws.Range('C40').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.999'
$ws.Range('E40').Value = '  -0.06%  '
$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.06'
$ws.Range('E41').Value = '  +11.06%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '133.93'
$ws.Range('E42').Value = '  +2.97%  '
$ws.Range('E43').Value = '  -0.05%  '
$ws.Range('E44').Value = '  +1.93%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.94'
$ws.Range('E45').Value = '  -4.70%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '16.87'
$ws.Range('E46').Value = '  -3.10%  '
$ws.Range('E47').Value = '  -4.61%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '21.83'
$ws.Range('E48').Value = '  -2.45%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.146.45'
$ws.Range('E49').Value = '  +2.53%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.07'
$ws.Range('E50').Value = '  -0.35%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.02'
$ws.Range('E51').Value = '  +24.81%  '
